$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "modified" timestamp in B21
$ws.Range("B21").Value = "2023-08-17T11:35:01+00:00"

# The skos:altLabel values that were (incorrectly) placed in column C
# (skos:altLabel) are actually skos:broader values and belong in column F
# (skos:broader). Move C24:C37 down-and-over into F24:F37, clearing out
# the old C cells that held data (C25:C30, C32:C37).

$ws.Range("F24").Value = "vocab.1000"
$ws.Range("F25").Value = "vocab.1000"
$ws.Range("F26").Value = "vocab.1000"
$ws.Range("F27").Value = "vocab.1000,vocab.1007"
$ws.Range("F28").Value = "vocab.1000,vocab.1007"
$ws.Range("F29").Value = "vocab.1000,vocab.1007"
$ws.Range("F31").Value = "vocab.1007"
$ws.Range("F32").Value = "vocab.1007"
$ws.Range("F33").Value = "vocab.1007"
$ws.Range("F34").Value = "vocab.1007"
$ws.Range("F35").Value = "vocab.1007"
$ws.Range("F36").Value = "vocab.1007"

$ws.Range("C25").Value = ""
$ws.Range("C26").Value = ""
$ws.Range("C27").Value = ""
$ws.Range("C28").Value = ""
$ws.Range("C29").Value = ""
$ws.Range("C30").Value = ""
$ws.Range("C32").Value = ""
$ws.Range("C33").Value = ""
$ws.Range("C34").Value = ""
$ws.Range("C35").Value = ""
$ws.Range("C36").Value = ""
$ws.Range("C37").Value = ""
